# Renames several label cells (shared strings) to shorter / cleaner
# descriptions ("Calculate Adj Net Income and add to dcf_inputs" cleanup
# pass), and leaves the cursor positioned on each sheet where the author's
# edit ended up.

$wb = $excel.ActiveWorkbook

$wsInputs = $wb.Worksheets.Item("Inputs")
$wsBalance = $wb.Worksheets.Item("Balance Sheet")
$wsIncome = $wb.Worksheets.Item("Income Statement")
$wsCashFlow = $wb.Worksheets.Item("Cash Flow Statement")

# --- Income Statement sheet --------------------------------------------
$wsIncome.Range("A3").Value = "Interest Income"
$wsIncome.Activate()
$wsIncome.Range("A13").Select()

# --- Balance Sheet sheet ------------------------------------------------
$wsBalance.Range("A4").Value = "Shares Outstanding"
$wsBalance.Activate()
$wsBalance.Range("A13").Select()

# --- Cash Flow Statement sheet -----------------------------------------
$wsCashFlow.Range("A3").Value = "Depreciation"
$wsCashFlow.Range("A5").Value = "Net New Debt"
$wsCashFlow.Range("A6").Value = "Net Stock Issuance"
$wsCashFlow.Activate()
$wsCashFlow.Range("A7").Select()

# --- Inputs sheet --------------------------------------------------------
$wsInputs.Range("A7").Value = "Stable Period Growth Rate"
$wsInputs.Activate()
$wsInputs.Range("A16").Select()
